$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay text (matching the source workbook, where every data
# cell is stored as a string) by pre-setting a text number format.
$textCells = @("D5", "D6", "D10", "D11", "D14", "D19", "D21", "D23", "D25", "D26", "D29", "D30", "D31", "D32", "D34", "D37", "D40", "D42", "D44", "D46", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '65.653.36'
$ws.Range("E2").Value = '  -0.14%  '
$ws.Range("D3").Value = '2.650.89'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '597.62'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").Value = '156.49'
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '0.397'
$ws.Range("E10").Value = '  -0.90%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '5.82'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").Value = '28.64'
$ws.Range("E14").Value = '  -2.83%  '
$ws.Range("D15").Value = '3.126.45'
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").Value = '65.475.60'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '2.660.99'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("D21").Value = '348.58'
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D23").Value = '68.98'
$ws.Range("E24").Value = '  +2.62%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '9.63'
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = '1.69'
$ws.Range("E26").Value = '  +3.93%  '
$ws.Range("E27").Value = '  -0.78%  '
$ws.Range("E28").Value = '  -2.32%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '7.92'
$ws.Range("E30").Value = '  -2.62%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.12'
$ws.Range("E31").Value = '  -1.22%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '530.40'
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '6.41'
$ws.Range("E34").Value = '  -2.52%  '
$ws.Range("E35").Value = '  -1.12%  '
$ws.Range("E36").Value = '  -1.15%  '
$ws.Range("D37").Value = '20.36'
$ws.Range("E37").Value = '  -0.29%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").Value = '155.25'
$ws.Range("E40").Value = '  -2.09%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '161.12'
$ws.Range("E42").Value = '  -2.26%  '
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").Value = '0.0606'
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("E45").Value = '  -1.15%  '
$ws.Range("D46").Value = '22.57'
$ws.Range("E46").Value = '  -2.78%  '
$ws.Range("E47").Value = '  -1.59%  '
$ws.Range("E48").Value = '  -1.80%  '
$ws.Range("D49").Value = '0.0994'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("E50").Value = '  +7.42%  '
$ws.Range("D51").Value = '19.75'
$ws.Range("E51").Value = '  -2.56%  '
